# Test.xlsx bugfix: image extensions png -> jpg + add "Links|webaddr" column
# with the corresponding web addresses (as real hyperlinks) on the "Cities"
# sheet, and reset the "Citizens" sheet scroll position back to column A.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Citizens sheet: just scroll the frozen/top-left cell back to A1.
# ---------------------------------------------------------------------
$wsCitizens = $wb.Worksheets.Item("Citizens")
$wsCitizens.Application.ActiveWindow.ScrollRow = 1
$wsCitizens.Application.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------
# Cities sheet: fix the picture filename extensions + add the new
# "Links|webaddr" hyperlink column.
# ---------------------------------------------------------------------
$wsCities = $wb.Worksheets.Item("Cities")

# Picture file name extensions: .png -> .jpg
$wsCities.Range("C2").Value = "City1.jpg"
$wsCities.Range("C3").Value = "City2.jpg|River.jpg"
$wsCities.Range("C4").Value = "River.jpg"

# New header cell, styled like the existing bold header row.
$wsCities.Range("D1").Value = "Links|webaddr"
$wsCities.Range("D1").Font.Bold = $true
$wsCities.Range("D1").Font.Size = 14

# New hyperlink cells. Adding a hyperlink auto-applies the blue/underline
# "Hyperlink" style, which the source file does not use here (the cells
# keep the plain default look) -- so copy the plain formatting from a
# normal cell (A2, style 0) back onto each new cell afterwards.
$wsCities.Hyperlinks.Add($wsCities.Range("D2"), "http://commitment.cornell.edu/", [Type]::Missing, [Type]::Missing, "http://commitment.cornell.edu/") | Out-Null
$wsCities.Hyperlinks.Add($wsCities.Range("D3"), "http://www.presidence.cg/accueil/", [Type]::Missing, [Type]::Missing, "http://www.presidence.cg/accueil/") | Out-Null
$wsCities.Hyperlinks.Add($wsCities.Range("D4"), "http://www.sydneydance.com.au/", [Type]::Missing, [Type]::Missing, "http://www.sydneydance.com.au/") | Out-Null

$wsCities.Range("A2").Copy() | Out-Null
$wsCities.Range("D2:D4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active selection to D10, matching the saved view state.
$wsCities.Activate()
$wsCities.Range("D10").Select() | Out-Null
